$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) values between row 2 and row 4
$ws.Range("D2").Value = 44981
$ws.Range("D4").Value = 44980

# Swap the Volumen values between row 2 and row 4
$ws.Range("M2").Value = 30
$ws.Range("M4").Value = 50
